$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2023-06-13"
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 17533
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región del Maule'
$ws.Range("S2").Value = 974
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = "2023-06-13"
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 130
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14462
$ws.Range("R3").Value = 'Región del Maule'
$ws.Range("S3").Value = 803

$ws.Range("D4").Value = "2022-06-02"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 806

$ws.Range("D5").Value = "2022-06-07"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14400
$ws.Range("R5").Value = 'Región del Maule'
$ws.Range("S5").Value = 800

$ws.Range("D6").Value = "2021-05-14"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = "2023-06-15"
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 140
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 19000
$ws.Range("P7").Value = 18429
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 1024
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = "2021-05-18"
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11500
$ws.Range("Q8").Value = '$/caja 12 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 11500
$ws.Range("T8").Value = 1

$ws.Range("D9").Value = "2022-06-30"
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 806

$ws.Range("D10").Value = "2022-05-27"
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12571
$ws.Range("Q10").Value = '$/caja 12 kilos empedrada'
$ws.Range("S10").Value = 1048
$ws.Range("T10").Value = 12

$ws.Range("D11").Value = "2022-05-26"
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 12500
$ws.Range("Q11").Value = '$/caja 12 kilos empedrada'
$ws.Range("S11").Value = 1042
$ws.Range("T11").Value = 12

$ws.Range("D12").Value = "2023-06-16"
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17429
$ws.Range("S12").Value = 968

$ws.Range("D13").Value = "2023-06-07"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17500
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Región del Maule'
$ws.Range("S13").Value = 972
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = "2023-06-14"
$ws.Range("M14").Value = 220
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 19000
$ws.Range("P14").Value = 18455
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = 'Provincia de Curicó'
$ws.Range("S14").Value = 1025
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = "2023-06-14"
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 833

$ws.Range("D16").Value = "2023-06-23"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1000

$ws.Range("D17").Value = "2023-06-20"
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 18000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19000
$ws.Range("S17").Value = 1056

$ws.Range("D18").Value = "2023-07-07"
$ws.Range("N18").Value = 18000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 18000
$ws.Range("S18").Value = 1000

$ws.Range("D19").Value = "2023-06-29"
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 80
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 17625
$ws.Range("R19").Value = 'Región del Maule'
$ws.Range("S19").Value = 979

$ws.Range("D20").Value = "2023-05-31"
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 12857
$ws.Range("Q20").Value = '$/caja 12 kilos granel'
$ws.Range("S20").Value = 12857
$ws.Range("T20").Value = 1

$ws.Range("D21").Value = "2023-05-31"
$ws.Range("L21").Value = 'Segunda'
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 11000
$ws.Range("Q21").Value = '$/caja 12 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 11000
$ws.Range("T21").Value = 1
